# The paragraph currently reads (in one of its colored runs):
#   "... Version for branch alternate 1)"
# The run holding " 1" (a leading space followed by the digit "1") must be
# split into two runs: one holding just " " (keeping the original run's
# identity) and a new run holding "2" (replacing "1"), while the following
# ")" run must stay exactly as it was.
#
# We locate " 1" via Find (read-only, no replacement) to get a Range with
# the correct character offsets, then operate on the single character "1"
# so the existing run splits cleanly around it.

$d = $word.ActiveDocument

$findRange = $d.Content.Duplicate
$findRange.Find.Execute(" 1)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)

$digitStart = $findRange.Start + 1
$digitEnd = $digitStart + 1

# Range covering just the "1" character.
$digit = $d.Range($digitStart, $digitEnd)

# Briefly toggling a character property on just the "1" forces Word to
# split its run away from the leading-space character that precedes it
# (the space keeps the original run, "1" becomes its own run).
$digit.Font.Bold = 1
$d.Range($digitStart, $digitEnd).Font.Bold = 0

# Replace "1" with "2" while the character is still its own run and force
# a split away from the ")" that follows, so the new "2" run does not
# re-merge with the trailing run.
$digit2 = $d.Range($digitStart, $digitEnd)
$digit2.Font.Bold = 1
$digit2.Text = "2"
$d.Range($digitStart, $digitEnd).Font.Bold = 0
